$wb = $excel.ActiveWorkbook

# Select a cell on the "mets" sheet (no data change there, just a view/selection change)
$metsWs = $wb.Worksheets.Item("mets")
$metsWs.Activate()
$metsWs.Range("D23").Select()

# Remove the "m4" and "m22" rows from the metsData sheet (input validation test fixture update)
$ws = $wb.Worksheets.Item("metsData")
$ws.Activate()
$ws.Rows.Item(5).Delete()
$ws.Rows.Item(22).Delete()
$ws.Range("A22").Select()
